$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2X project definition")

# Fix typo in the "GPS Enable" feature description (row 34, column D)
$ws.Range("D34").Value = "The GPS system can be started and stopped by AT cmd or SIM USB"

# Mark additional rows as "Developed" (column E) per updated documentation status
$ws.Range("E15").Value = "X"
$ws.Range("E16").Value = "X"
$ws.Range("E17").Value = "X"
$ws.Range("E18").Value = "X"
$ws.Range("E21").Value = "X"
$ws.Range("E22").Value = "p"
$ws.Range("E23").Value = "X"
$ws.Range("E26").Value = "p"
$ws.Range("E29").Value = "X"
$ws.Range("E33").Value = "X"
$ws.Range("E50").Value = "X"
$ws.Range("E51").Value = "X"
$ws.Range("E58").Value = "X"
$ws.Range("E59").Value = "X"

# Update the active selection to reflect where editing ended
$ws.Range("E63").Select()
